$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 1: "Odd_CS_4-4" moves from AM1 to AG1, shifting the old
#     AG1:AL1 block (Odd_CS_0-1 .. Odd_CS_2-3) one column to the right
#     (into AH1:AM1). Columns AG..AM are columns 33..39. ---
$headerValues = @(
    "Odd_CS_4-4",
    "Odd_CS_0-1",
    "Odd_CS_0-2",
    "Odd_CS_1-2",
    "Odd_CS_0-3",
    "Odd_CS_1-3",
    "Odd_CS_2-3"
)
$startCol = 33  # AG
for ($i = 0; $i -lt $headerValues.Length; $i++) {
    $ws.Cells.Item(1, $startCol + $i).Value = $headerValues[$i]
}

# --- Row 2: replace with the new match data (columns A..BD = 1..56) ---
$row2 = @(
    "xCf4Akf2", "29/10/2024", "07:45", "MALAYSIA - SUPER LEAGUE", "Sabah", "Kedah",
    1.5, 3.85, 5.5, 2, 2.25, 5.3, 1.03, 11.9, 1.16, 3.86, 1.65, 2,
    1.31, 3.26, 1.76, 2.01, 6.3, 6.4, 6.8, 9, 9.75, 18.5, 12, 6.8, 13,
    50, 300, 14.5, 30, 14.5, 80, 40, 37, 3.4, 7.1, 15.5, 21, 45, 200,
    3, 7.3, 60, 7.3, 30, 32, 200, 200, 350, 51, 51
)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

# --- Row 3: delete entirely (rows below shift up; here it just vanishes) ---
$ws.Rows.Item(3).Delete()
